$wb = $excel.ActiveWorkbook

# facil_elo_p3_c3 (sheet index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = -0.064
$ws.Range("I2").Value = 0.949
$ws.Range("H3").Value = 0.8110000000000001
$ws.Range("I3").Value = 0.417
$ws.Range("H4").Value = 0.066
$ws.Range("I4").Value = 0.947
$ws.Range("H5").Value = 0.42
$ws.Range("I5").Value = 0.675
$ws.Range("H6").Value = 0.162
$ws.Range("I6").Value = 0.871
$ws.Range("H7").Value = -0.289
$ws.Range("I7").Value = 0.772
$ws.Range("H8").Value = 0.288
$ws.Range("I8").Value = 0.774
$ws.Range("H9").Value = -1.06
$ws.Range("I9").Value = 0.289
$ws.Range("H10").Value = 0.572
$ws.Range("I10").Value = 0.5669999999999999
$ws.Range("H11").Value = 0.122
$ws.Range("I11").Value = 0.903
$ws.Range("H12").Value = -0.058
$ws.Range("I12").Value = 0.953
$ws.Range("H13").Value = 0.223
$ws.Range("I13").Value = 0.824
$ws.Range("H14").Value = -2.028
$ws.Range("I14").Value = 0.043
$ws.Range("H15").Value = -0.512
$ws.Range("I15").Value = 0.609
$ws.Range("H16").Value = 1.168
$ws.Range("I16").Value = 0.243
$ws.Range("H17").Value = 1.115
$ws.Range("I17").Value = 0.265
$ws.Range("H18").Value = 1.684
$ws.Range("I18").Value = 0.092
$ws.Range("H19").Value = -0.727
$ws.Range("I19").Value = 0.467

# facil_elo_p4_c5 (sheet index 10)
$ws = $wb.Worksheets.Item(10)
$ws.Range("I2").Value = -1.529
$ws.Range("J2").Value = 0.126
$ws.Range("I3").Value = 0.949
$ws.Range("J3").Value = 0.342
$ws.Range("I4").Value = -0.054
$ws.Range("J4").Value = 0.957
$ws.Range("I5").Value = 0.246
$ws.Range("J5").Value = 0.805
$ws.Range("I6").Value = 2.096
$ws.Range("J6").Value = 0.036
$ws.Range("I7").Value = 0.282
$ws.Range("J7").Value = 0.778
$ws.Range("I8").Value = 2.119
$ws.Range("J8").Value = 0.034
$ws.Range("I9").Value = -0.324
$ws.Range("J9").Value = 0.746
$ws.Range("I10").Value = 0.882
$ws.Range("J10").Value = 0.378
$ws.Range("I11").Value = -0.097
$ws.Range("J11").Value = 0.923
$ws.Range("I12").Value = -1.81
$ws.Range("J12").Value = 0.07000000000000001
$ws.Range("I13").Value = -0.438
$ws.Range("J13").Value = 0.662
$ws.Range("I14").Value = -1.801
$ws.Range("J14").Value = 0.07199999999999999
$ws.Range("I15").Value = -0.448
$ws.Range("J15").Value = 0.654
$ws.Range("I16").Value = -0.362
$ws.Range("J16").Value = 0.718
$ws.Range("I17").Value = 0.987
$ws.Range("J17").Value = 0.324
$ws.Range("I18").Value = 1.313
$ws.Range("J18").Value = 0.189
$ws.Range("I19").Value = 0.156
$ws.Range("J19").Value = 0.876

# facil_elo_p5_c5 (sheet index 11)
$ws = $wb.Worksheets.Item(11)
$ws.Range("J2").Value = -1.682
$ws.Range("K2").Value = 0.093
$ws.Range("J3").Value = 0.162
$ws.Range("K3").Value = 0.871
$ws.Range("J4").Value = 0.015
$ws.Range("K4").Value = 0.988
$ws.Range("J5").Value = 0.091
$ws.Range("K5").Value = 0.927
$ws.Range("J6").Value = 1.447
$ws.Range("K6").Value = 0.148
$ws.Range("J7").Value = 0.023
$ws.Range("K7").Value = 0.982
$ws.Range("J8").Value = 2.124
$ws.Range("K8").Value = 0.034
$ws.Range("J9").Value = 0.219
$ws.Range("K9").Value = 0.826
$ws.Range("J10").Value = 0.602
$ws.Range("K10").Value = 0.547
$ws.Range("J11").Value = 0.081
$ws.Range("K11").Value = 0.9360000000000001
$ws.Range("J12").Value = -1.156
$ws.Range("K12").Value = 0.248
$ws.Range("J13").Value = -0.243
$ws.Range("K13").Value = 0.8080000000000001
$ws.Range("J14").Value = -1.571
$ws.Range("K14").Value = 0.116
$ws.Range("J15").Value = -0.288
$ws.Range("K15").Value = 0.773
$ws.Range("J17").Value = 0.9
$ws.Range("K17").Value = 0.368
$ws.Range("J18").Value = 0.694
$ws.Range("K18").Value = 0.488
$ws.Range("J19").Value = -0.493
$ws.Range("K19").Value = 0.622

# facil_elo_p6_c5 (sheet index 12)
$ws = $wb.Worksheets.Item(12)
$ws.Range("J2").Value = -1.053
$ws.Range("K2").Value = 0.293
$ws.Range("J3").Value = 0.13
$ws.Range("K3").Value = 0.896
$ws.Range("J4").Value = -0.399
$ws.Range("K4").Value = 0.6899999999999999
$ws.Range("J5").Value = -0.323
$ws.Range("K5").Value = 0.746
$ws.Range("J6").Value = 1.047
$ws.Range("K6").Value = 0.295
$ws.Range("J7").Value = -0.004
$ws.Range("K7").Value = 0.996
$ws.Range("J8").Value = 1.03
$ws.Range("K8").Value = 0.303
$ws.Range("J9").Value = 0.056
$ws.Range("K9").Value = 0.955
$ws.Range("J10").Value = 0.854
$ws.Range("K10").Value = 0.393
$ws.Range("J11").Value = 0.648
$ws.Range("K11").Value = 0.517
$ws.Range("J12").Value = -1.021
$ws.Range("K12").Value = 0.307
$ws.Range("J13").Value = -0.106
$ws.Range("K13").Value = 0.916
$ws.Range("J14").Value = -1.552
$ws.Range("K14").Value = 0.121
$ws.Range("J15").Value = -0.076
$ws.Range("K15").Value = 0.9389999999999999
$ws.Range("J16").Value = -0.411
$ws.Range("K16").Value = 0.681
$ws.Range("J17").Value = 0.723
$ws.Range("K17").Value = 0.47
$ws.Range("J18").Value = 0.074
$ws.Range("K18").Value = 0.9409999999999999
$ws.Range("J19").Value = -0.738
$ws.Range("K19").Value = 0.46

# facil_elo_p3_c6 (sheet index 13)
$ws = $wb.Worksheets.Item(13)
$ws.Range("I2").Value = -2.282
$ws.Range("J2").Value = 0.022
$ws.Range("I3").Value = -0.437
$ws.Range("J3").Value = 0.662
$ws.Range("I4").Value = -0.746
$ws.Range("J4").Value = 0.455
$ws.Range("I5").Value = 1.511
$ws.Range("J5").Value = 0.131
$ws.Range("I6").Value = 1.251
$ws.Range("J6").Value = 0.211
$ws.Range("I7").Value = 0.322
$ws.Range("J7").Value = 0.747
$ws.Range("I8").Value = 1.995
$ws.Range("J8").Value = 0.046
$ws.Range("I9").Value = 0.967
$ws.Range("J9").Value = 0.334
$ws.Range("I10").Value = 1.344
$ws.Range("J10").Value = 0.179
$ws.Range("I11").Value = -1.444
$ws.Range("J11").Value = 0.149
$ws.Range("I12").Value = -1.753
$ws.Range("J12").Value = 0.08
$ws.Range("I13").Value = -0.509
$ws.Range("J13").Value = 0.61
$ws.Range("I14").Value = -2.584
$ws.Range("J14").Value = 0.01
$ws.Range("I15").Value = -0.575
$ws.Range("J15").Value = 0.5649999999999999
$ws.Range("I16").Value = -0.448
$ws.Range("J16").Value = 0.654
$ws.Range("I17").Value = -0.055
$ws.Range("J17").Value = 0.956
$ws.Range("I18").Value = 2.479
$ws.Range("J18").Value = 0.013
$ws.Range("I19").Value = 0.362
$ws.Range("J19").Value = 0.717

# facil_elo_p4_c6 (sheet index 14)
$ws = $wb.Worksheets.Item(14)
$ws.Range("I2").Value = -2.246
$ws.Range("J2").Value = 0.025
$ws.Range("I3").Value = -0.251
$ws.Range("J3").Value = 0.802
$ws.Range("I4").Value = -1.216
$ws.Range("J4").Value = 0.224
$ws.Range("I5").Value = 0.665
$ws.Range("J5").Value = 0.506
$ws.Range("I6").Value = 1.93
$ws.Range("J6").Value = 0.054
$ws.Range("I7").Value = 0.953
$ws.Range("J7").Value = 0.341
$ws.Range("I8").Value = 2.155
$ws.Range("J8").Value = 0.031
$ws.Range("I9").Value = 0.739
$ws.Range("J9").Value = 0.46
$ws.Range("I10").Value = 1.672
$ws.Range("J10").Value = 0.094
$ws.Range("I11").Value = -0.573
$ws.Range("J11").Value = 0.5669999999999999
$ws.Range("I12").Value = -2.379
$ws.Range("J12").Value = 0.017
$ws.Range("I13").Value = -1.157
$ws.Range("J13").Value = 0.247
$ws.Range("I14").Value = -1.953
$ws.Range("J14").Value = 0.051
$ws.Range("I15").Value = 0.063
$ws.Range("J15").Value = 0.95
$ws.Range("I16").Value = -0.14
$ws.Range("J16").Value = 0.888
$ws.Range("I17").Value = 0.964
$ws.Range("J17").Value = 0.335
$ws.Range("I18").Value = 2.152
$ws.Range("J18").Value = 0.031
$ws.Range("I19").Value = -0.032
$ws.Range("J19").Value = 0.975

# facil_elo_p5_c6 (sheet index 15)
$ws = $wb.Worksheets.Item(15)
$ws.Range("J2").Value = -2.655
$ws.Range("K2").Value = 0.008
$ws.Range("J3").Value = -0.553
$ws.Range("K3").Value = 0.58
$ws.Range("J4").Value = -0.587
$ws.Range("K4").Value = 0.5570000000000001
$ws.Range("J5").Value = 0.59
$ws.Range("K5").Value = 0.555
$ws.Range("J6").Value = 1.743
$ws.Range("K6").Value = 0.081
$ws.Range("J7").Value = 0.64
$ws.Range("K7").Value = 0.522
$ws.Range("J8").Value = 2.73
$ws.Range("K8").Value = 0.006
$ws.Range("J9").Value = 0.8100000000000001
$ws.Range("K9").Value = 0.418
$ws.Range("J10").Value = 1.144
$ws.Range("K10").Value = 0.252
$ws.Range("J11").Value = -0.63
$ws.Range("K11").Value = 0.529
$ws.Range("J12").Value = -2.038
$ws.Range("K12").Value = 0.042
$ws.Range("J13").Value = -0.759
$ws.Range("K13").Value = 0.448
$ws.Range("J14").Value = -1.541
$ws.Range("K14").Value = 0.123
$ws.Range("J15").Value = -0.015
$ws.Range("K15").Value = 0.988
$ws.Range("J16").Value = -0.297
$ws.Range("K16").Value = 0.767
$ws.Range("J17").Value = 0.947
$ws.Range("K17").Value = 0.344
$ws.Range("J18").Value = 1.376
$ws.Range("K18").Value = 0.169
$ws.Range("J19").Value = -0.5600000000000001
$ws.Range("K19").Value = 0.576

# facil_elo_p6_c6 (sheet index 16)
$ws = $wb.Worksheets.Item(16)
$ws.Range("J2").Value = -1.898
$ws.Range("K2").Value = 0.058
$ws.Range("J3").Value = -0.947
$ws.Range("K3").Value = 0.344
$ws.Range("J4").Value = -1.663
$ws.Range("K4").Value = 0.096
$ws.Range("J5").Value = 0.181
$ws.Range("K5").Value = 0.857
$ws.Range("J6").Value = 2.076
$ws.Range("K6").Value = 0.038
$ws.Range("J7").Value = 0.591
$ws.Range("K7").Value = 0.555
$ws.Range("J8").Value = 2.071
$ws.Range("K8").Value = 0.038
$ws.Range("J9").Value = 0.955
$ws.Range("K9").Value = 0.34
$ws.Range("J10").Value = 2.135
$ws.Range("K10").Value = 0.033
$ws.Range("J11").Value = 0.154
$ws.Range("K11").Value = 0.877
$ws.Range("J12").Value = -2.68
$ws.Range("K12").Value = 0.007
$ws.Range("J13").Value = -0.657
$ws.Range("K13").Value = 0.511
$ws.Range("J14").Value = -1.36
$ws.Range("K14").Value = 0.174
$ws.Range("J15").Value = 0.251
$ws.Range("K15").Value = 0.802
$ws.Range("J16").Value = -0.803
$ws.Range("K16").Value = 0.422
$ws.Range("J17").Value = 0.779
$ws.Range("K17").Value = 0.436
$ws.Range("J18").Value = 0.856
$ws.Range("K18").Value = 0.392
$ws.Range("J19").Value = -0.537
$ws.Range("K19").Value = 0.592

# facil_elo_p4_c3 (sheet index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = -0.35
$ws.Range("I2").Value = 0.726
$ws.Range("H3").Value = 0.508
$ws.Range("I3").Value = 0.611
$ws.Range("H4").Value = -1.064
$ws.Range("I4").Value = 0.287
$ws.Range("H5").Value = 0.465
$ws.Range("I5").Value = 0.642
$ws.Range("H6").Value = 1.249
$ws.Range("I6").Value = 0.212
$ws.Range("H7").Value = 0.243
$ws.Range("I7").Value = 0.8080000000000001
$ws.Range("H8").Value = 0.702
$ws.Range("I8").Value = 0.483
$ws.Range("H9").Value = -0.328
$ws.Range("I9").Value = 0.743
$ws.Range("H10").Value = 1.369
$ws.Range("I10").Value = 0.171
$ws.Range("H11").Value = 0.209
$ws.Range("I11").Value = 0.835
$ws.Range("H12").Value = -0.96
$ws.Range("I12").Value = 0.337
$ws.Range("H13").Value = -0.371
$ws.Range("I13").Value = 0.711
$ws.Range("H14").Value = -2.069
$ws.Range("I14").Value = 0.039
$ws.Range("H15").Value = -0.438
$ws.Range("I15").Value = 0.661
$ws.Range("H16").Value = 0.601
$ws.Range("I16").Value = 0.548
$ws.Range("H17").Value = 1.299
$ws.Range("I17").Value = 0.194
$ws.Range("H18").Value = 2.421
$ws.Range("I18").Value = 0.015
$ws.Range("H19").Value = -0.6840000000000001
$ws.Range("I19").Value = 0.494

# facil_elo_p5_c3 (sheet index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H2").Value = -0.094
$ws.Range("I2").Value = 0.925
$ws.Range("H3").Value = -0.153
$ws.Range("I3").Value = 0.878
$ws.Range("H4").Value = -0.838
$ws.Range("I4").Value = 0.402
$ws.Range("H5").Value = -0.161
$ws.Range("I5").Value = 0.872
$ws.Range("I6").Value = 0.374
$ws.Range("H7").Value = 0.242
$ws.Range("I7").Value = 0.8090000000000001
$ws.Range("H8").Value = 0.336
$ws.Range("I8").Value = 0.737
$ws.Range("H9").Value = 0.31
$ws.Range("I9").Value = 0.756
$ws.Range("H10").Value = 1.228
$ws.Range("I10").Value = 0.219
$ws.Range("H11").Value = 0.707
$ws.Range("I11").Value = 0.48
$ws.Range("H12").Value = -0.529
$ws.Range("I12").Value = 0.597
$ws.Range("H13").Value = -0.379
$ws.Range("I13").Value = 0.705
$ws.Range("H14").Value = -1.012
$ws.Range("I14").Value = 0.312
$ws.Range("H15").Value = -0.157
$ws.Range("I15").Value = 0.875
$ws.Range("H16").Value = 0.85
$ws.Range("I16").Value = 0.395
$ws.Range("H17").Value = 1.435
$ws.Range("I17").Value = 0.151
$ws.Range("H18").Value = 1.43
$ws.Range("I18").Value = 0.153
$ws.Range("H19").Value = -0.977
$ws.Range("I19").Value = 0.329

# facil_elo_p6_c3 (sheet index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("I2").Value = -0.364
$ws.Range("J2").Value = 0.716
$ws.Range("J3").Value = 0.674
$ws.Range("J4").Value = 0.367
$ws.Range("I5").Value = -0.596
$ws.Range("J5").Value = 0.552
$ws.Range("I6").Value = 0.642
$ws.Range("J6").Value = 0.521
$ws.Range("I7").Value = 0.01
$ws.Range("J7").Value = 0.992
$ws.Range("I8").Value = 0.191
$ws.Range("J8").Value = 0.848
$ws.Range("I9").Value = 0.522
$ws.Range("J9").Value = 0.602
$ws.Range("I10").Value = 1.035
$ws.Range("J10").Value = 0.301
$ws.Range("I11").Value = 1.286
$ws.Range("J11").Value = 0.198
$ws.Range("I12").Value = -0.44
$ws.Range("J12").Value = 0.66
$ws.Range("I13").Value = 0.045
$ws.Range("J13").Value = 0.964
$ws.Range("I14").Value = -1.198
$ws.Range("J14").Value = 0.231
$ws.Range("I15").Value = -0.24
$ws.Range("J15").Value = 0.8100000000000001
$ws.Range("I16").Value = -0.473
$ws.Range("J16").Value = 0.636
$ws.Range("I17").Value = 0.705
$ws.Range("J17").Value = 0.481
$ws.Range("I18").Value = 1.066
$ws.Range("J18").Value = 0.287
$ws.Range("I19").Value = -1.112
$ws.Range("J19").Value = 0.266

# facil_elo_p3_c4 (sheet index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = -1.299
$ws.Range("I2").Value = 0.194
$ws.Range("H3").Value = 0.889
$ws.Range("I3").Value = 0.374
$ws.Range("H4").Value = -0.848
$ws.Range("I4").Value = 0.396
$ws.Range("H5").Value = -0.357
$ws.Range("I5").Value = 0.721
$ws.Range("H6").Value = 0.141
$ws.Range("I6").Value = 0.888
$ws.Range("H7").Value = -1.208
$ws.Range("I7").Value = 0.227
$ws.Range("H8").Value = 1.285
$ws.Range("I8").Value = 0.199
$ws.Range("H9").Value = -0.837
$ws.Range("I9").Value = 0.402
$ws.Range("H10").Value = 1.261
$ws.Range("I10").Value = 0.207
$ws.Range("H11").Value = 0.595
$ws.Range("I11").Value = 0.552
$ws.Range("H12").Value = -0.194
$ws.Range("I12").Value = 0.846
$ws.Range("H13").Value = 0.967
$ws.Range("I13").Value = 0.334
$ws.Range("H14").Value = -1.962
$ws.Range("I14").Value = 0.05
$ws.Range("H15").Value = -1.19
$ws.Range("I15").Value = 0.234
$ws.Range("H16").Value = 0.657
$ws.Range("I16").Value = 0.511
$ws.Range("H17").Value = 0.296
$ws.Range("I17").Value = 0.767
$ws.Range("H18").Value = 1.358
$ws.Range("I18").Value = 0.174
$ws.Range("H19").Value = 1.118
$ws.Range("I19").Value = 0.263

# facil_elo_p4_c4 (sheet index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Range("I2").Value = -1.529
$ws.Range("J2").Value = 0.126
$ws.Range("I3").Value = 0.43
$ws.Range("J3").Value = 0.667
$ws.Range("I4").Value = -0.79
$ws.Range("J4").Value = 0.43
$ws.Range("I5").Value = -0.173
$ws.Range("J5").Value = 0.863
$ws.Range("I6").Value = 1.438
$ws.Range("J6").Value = 0.15
$ws.Range("I7").Value = -0.678
$ws.Range("J7").Value = 0.498
$ws.Range("I8").Value = 1.666
$ws.Range("J8").Value = 0.096
$ws.Range("I9").Value = 0.185
$ws.Range("J9").Value = 0.853
$ws.Range("I10").Value = 1.438
$ws.Range("J10").Value = 0.15
$ws.Range("I11").Value = 0.549
$ws.Range("J11").Value = 0.583
$ws.Range("I12").Value = -1.319
$ws.Range("J12").Value = 0.187
$ws.Range("I13").Value = 0.443
$ws.Range("J13").Value = 0.658
$ws.Range("I14").Value = -2.046
$ws.Range("J14").Value = 0.041
$ws.Range("I15").Value = -1.257
$ws.Range("J15").Value = 0.209
$ws.Range("I16").Value = 0.391
$ws.Range("J16").Value = 0.696
$ws.Range("I17").Value = 0.958
$ws.Range("J17").Value = 0.338
$ws.Range("I18").Value = 1.687
$ws.Range("J18").Value = 0.092
$ws.Range("I19").Value = 0.656
$ws.Range("J19").Value = 0.512

# facil_elo_p5_c4 (sheet index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Range("I2").Value = -1.44
$ws.Range("J2").Value = 0.15
$ws.Range("I3").Value = -0.077
$ws.Range("J3").Value = 0.9389999999999999
$ws.Range("I4").Value = -0.799
$ws.Range("J4").Value = 0.424
$ws.Range("I5").Value = -0.665
$ws.Range("J5").Value = 0.506
$ws.Range("I6").Value = 0.829
$ws.Range("J6").Value = 0.407
$ws.Range("I7").Value = -0.754
$ws.Range("J7").Value = 0.451
$ws.Range("I8").Value = 1.578
$ws.Range("J8").Value = 0.115
$ws.Range("I9").Value = 0.442
$ws.Range("J9").Value = 0.658
$ws.Range("I10").Value = 1.182
$ws.Range("J10").Value = 0.237
$ws.Range("I11").Value = 1.039
$ws.Range("J11").Value = 0.299
$ws.Range("I12").Value = -0.648
$ws.Range("J12").Value = 0.517
$ws.Range("I13").Value = 0.512
$ws.Range("J13").Value = 0.609
$ws.Range("I14").Value = -1.75
$ws.Range("J14").Value = 0.08
$ws.Range("I15").Value = -1.043
$ws.Range("J15").Value = 0.297
$ws.Range("I16").Value = 0.722
$ws.Range("J16").Value = 0.47
$ws.Range("I17").Value = 0.9
$ws.Range("J17").Value = 0.368
$ws.Range("I18").Value = 0.998
$ws.Range("J18").Value = 0.318
$ws.Range("I19").Value = 0.065
$ws.Range("J19").Value = 0.948

# facil_elo_p6_c4 (sheet index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Range("J2").Value = -1.168
$ws.Range("K2").Value = 0.243
$ws.Range("J3").Value = -0.012
$ws.Range("K3").Value = 0.991
$ws.Range("J4").Value = -1.184
$ws.Range("K4").Value = 0.236
$ws.Range("J5").Value = -1.085
$ws.Range("K5").Value = 0.278
$ws.Range("J6").Value = 0.606
$ws.Range("K6").Value = 0.545
$ws.Range("J7").Value = -0.865
$ws.Range("K7").Value = 0.387
$ws.Range("J8").Value = 0.919
$ws.Range("K8").Value = 0.358
$ws.Range("J9").Value = 0.215
$ws.Range("K9").Value = 0.83
$ws.Range("J10").Value = 1.283
$ws.Range("K10").Value = 0.2
$ws.Range("J11").Value = 1.636
$ws.Range("K11").Value = 0.102
$ws.Range("J12").Value = -0.53
$ws.Range("K12").Value = 0.596
$ws.Range("J13").Value = 0.704
$ws.Range("K13").Value = 0.482
$ws.Range("J14").Value = -1.441
$ws.Range("K14").Value = 0.149
$ws.Range("J15").Value = -0.921
$ws.Range("K15").Value = 0.357
$ws.Range("J16").Value = 0.122
$ws.Range("K16").Value = 0.903
$ws.Range("J17").Value = 0.53
$ws.Range("K17").Value = 0.596
$ws.Range("J18").Value = 0.703
$ws.Range("K18").Value = 0.482
$ws.Range("J19").Value = -0.064
$ws.Range("K19").Value = 0.949

# facil_elo_p3_c5 (sheet index 9)
$ws = $wb.Worksheets.Item(9)
$ws.Range("I2").Value = -1.045
$ws.Range("J2").Value = 0.296
$ws.Range("I3").Value = 1.108
$ws.Range("J3").Value = 0.268
$ws.Range("I4").Value = 0.059
$ws.Range("J4").Value = 0.953
$ws.Range("I5").Value = 0.902
$ws.Range("J5").Value = 0.367
$ws.Range("I6").Value = 0.931
$ws.Range("J6").Value = 0.352
$ws.Range("I7").Value = 0.099
$ws.Range("J7").Value = 0.921
$ws.Range("I8").Value = 1.447
$ws.Range("J8").Value = 0.148
$ws.Range("I9").Value = -0.806
$ws.Range("J9").Value = 0.42
$ws.Range("I10").Value = 0.6820000000000001
$ws.Range("J10").Value = 0.495
$ws.Range("I11").Value = -0.726
$ws.Range("J11").Value = 0.468
$ws.Range("I12").Value = -0.984
$ws.Range("J12").Value = 0.325
$ws.Range("I13").Value = -0.203
$ws.Range("J13").Value = 0.839
$ws.Range("I14").Value = -1.709
$ws.Range("J14").Value = 0.08699999999999999
$ws.Range("I15").Value = -0.836
$ws.Range("J15").Value = 0.403
$ws.Range("I16").Value = 0.039
$ws.Range("J16").Value = 0.969
$ws.Range("I17").Value = 0.301
$ws.Range("J17").Value = 0.763
$ws.Range("I18").Value = 1.135
$ws.Range("J18").Value = 0.256
$ws.Range("I19").Value = 0.678
$ws.Range("J19").Value = 0.498

